# Applies two logically-separate edits found in the target diff:
#
# 1) Three tables (on slides 14, 15 and 16) get their table style switched
#    from the custom "Table_0" style ({5BE1EE9A-F694-400E-92F9-7FCAB5BCD114})
#    to PowerPoint's built-in "No Style, No Grid" table style
#    ({466681D4-EFF8-41CD-8AF7-2FCB3B3695E8}).
#
# 2) The presentation's theme colour palette is swapped: the deck currently
#    uses the "Integral / Red Violet" colour palette; the edit changes it to
#    the standard "Office" colour palette (this is what actually differs
#    between ppt/theme/theme1.xml and ppt/theme/theme2.xml - every other
#    part of the two themes, fonts and format scheme, is already identical).

$p = $ppt.ActivePresentation

# --- 1) Re-style the three tables -----------------------------------------
$noStyleNoGrid = "{466681D4-EFF8-41CD-8AF7-2FCB3B3695E8}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    $tbl = $slide.Shapes.Item(1).Table
    $tbl.ApplyStyle($noStyleNoGrid)
}

# --- 2) Swap the theme colour scheme (Integral/Red Violet -> Office) ------
# Colour scheme slot order exposed through ThemeColorScheme.Colors(i):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
#   9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
